# Cardiology data workbook update
# - Row 10 ("Distribution type") is replaced with Shapiro-Wilk style p-value /
#   normality-test results for the numeric predictor columns.
# - Row 11 ("Is numeric data skewed? Type") gains a value in the "ca" (N) column.
# - Row 12 ("Level of correlation (Highest)") gains a value in the "ca" (N) column.
# - Several columns are widened to fit the new, longer text.
# - The active cell selection is moved to N13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: replace distribution-type notes with normality test p-values ---
$ws.Range("B10").Value = "p-value = 0.0007914 -> Deviates From Normality"
$ws.Range("E10").Value = "p-value = 8.431e-07 -> Normal Distribution"
$ws.Range("F10").Value = "p-value = 4.548e-09 -> Normal Distribution"
$ws.Range("I10").Value = "p-value = 3.104e-07 -> Normal Distribution"
$ws.Range("J10").Value = "p-value = 5.608e-05 -> Normal Distribution"
$ws.Range("L10").Value = "p-value < 2.2e-16 -> Deviates From Normal Distribution"
$ws.Range("N10").Value = "p-value < 2.2e-16 -> Deviates From Normality"

# --- Row 11: add the "ca" skew figure (others stay the same) ---
$ws.Range("N11").Value = "1.186547, Positively"

# --- Row 12: add the "ca" correlation figure (others stay the same) ---
$ws.Range("N12").Value = "0.3120721 age"

# --- Column widths: widen to accommodate the new text ---
$ws.Columns.Item(2).ColumnWidth = 27.333333333333332
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Columns.Item(4).ColumnWidth = 16.166666666666668
$ws.Columns.Item(5).ColumnWidth = 27
$ws.Columns.Item(6).ColumnWidth = 25.333333333333332
$ws.Columns.Item(7).ColumnWidth = 16.833333333333332
$ws.Columns.Item(8).ColumnWidth = 14.666666666666666
$ws.Columns.Item(9).ColumnWidth = 22.5
$ws.Columns.Item(10).ColumnWidth = 23.666666666666668
$ws.Columns.Item(12).ColumnWidth = 25.333333333333332
$ws.Columns.Item(14).ColumnWidth = 20.166666666666668

# --- Selection moves to N13 ---
$ws.Range("N13").Select() | Out-Null
